$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) ---
# "Volume 30   Number  45" -> "Volume 30   Number  46"
$ws.Range("A8").Value2 = "Volume 30   Number  46"
# "Report Covering the Week  11/6/2023  Through  11/12/2023"
#   -> "Report Covering the Week  11/13/2023  Through  11/19/2023"
$ws.Range("C9").Value2 = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Crime statistics table updates (rows 15-30) ---

# Row 15 (Rape)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value2 = '0'
$ws.Range("E15").Value2 = -100
$ws.Range("F15").Value2 = 2
$ws.Range("H15").Value2 = -50
$ws.Range("J15").Value2 = 31
$ws.Range("K15").Value2 = 3.225806451612
$ws.Range("M15").Value2 = -5.882352941176
$ws.Range("N15").Value2 = -58.974358974359

# Row 16 (Robbery)
$ws.Range("C16").Value2 = 9
$ws.Range("D16").Value2 = 5
$ws.Range("E16").Value2 = 80
$ws.Range("F16").Value2 = 26
$ws.Range("G16").Value2 = 24
$ws.Range("H16").Value2 = 8.333333333333
$ws.Range("I16").Value2 = 246
$ws.Range("J16").Value2 = 307
$ws.Range("K16").Value2 = -19.869706840390
$ws.Range("L16").Value2 = 21.782178217821
$ws.Range("M16").Value2 = -43.577981651376
$ws.Range("N16").Value2 = -87.809712586719

# Row 17 (Fel. Assault)
$ws.Range("C17").Value2 = 13
$ws.Range("D17").Value2 = 12
$ws.Range("E17").Value2 = 8.333333333333
$ws.Range("F17").Value2 = 48
$ws.Range("G17").Value2 = 57
$ws.Range("H17").Value2 = -15.789473684210
$ws.Range("I17").Value2 = 615
$ws.Range("J17").Value2 = 665
$ws.Range("K17").Value2 = -7.518796992481
$ws.Range("L17").Value2 = 10.017889087656
$ws.Range("M17").Value2 = 17.590822179732
$ws.Range("N17").Value2 = -40.751445086705

# Row 18 (Burglary)
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = -33.333333333333
$ws.Range("F18").Value2 = 6
$ws.Range("G18").Value2 = 19
$ws.Range("H18").Value2 = -68.421052631578
$ws.Range("I18").Value2 = 113
$ws.Range("J18").Value2 = 177
$ws.Range("K18").Value2 = -36.158192090395
$ws.Range("L18").Value2 = -13.740458015267
$ws.Range("M18").Value2 = -54.251012145749
$ws.Range("N18").Value2 = -83.670520231213

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value2 = 3
$ws.Range("D19").Value2 = 9
$ws.Range("E19").Value2 = -66.666666666666
$ws.Range("F19").Value2 = 23
$ws.Range("H19").Value2 = -45.238095238095
$ws.Range("I19").Value2 = 298
$ws.Range("J19").Value2 = 364
$ws.Range("K19").Value2 = -18.131868131868
$ws.Range("L19").Value2 = 7.194244604316
$ws.Range("M19").Value2 = -25.5
$ws.Range("N19").Value2 = -65.023474178403

# Row 20 (G.L.A.)
$ws.Range("C20").Value2 = 6
$ws.Range("D20").Value2 = 1
$ws.Range("E20").Value2 = 500
$ws.Range("F20").Value2 = 13
$ws.Range("G20").Value2 = 14
$ws.Range("H20").Value2 = -7.142857142857
$ws.Range("I20").Value2 = 136
$ws.Range("J20").Value2 = 173
$ws.Range("K20").Value2 = -21.387283236994
$ws.Range("L20").Value2 = -7.482993197278
$ws.Range("M20").Value2 = 1.492537313432
$ws.Range("N20").Value2 = -78.205128205128

# Row 21 (TOTAL)
$ws.Range("C21").Value2 = 33
$ws.Range("D21").Value2 = 31
$ws.Range("E21").Value2 = 6.451612903225
$ws.Range("F21").Value2 = 118
$ws.Range("G21").Value2 = 162
$ws.Range("H21").Value2 = -27.160493827160
$ws.Range("I21").Value2 = 1446
$ws.Range("J21").Value2 = 1736
$ws.Range("K21").Value2 = -16.705069124424
$ws.Range("L21").Value2 = 5.934065934065
$ws.Range("M21").Value2 = -19.577308120133
$ws.Range("N21").Value2 = -73.062593144560

# Row 22 (Transit)
$ws.Range("C22").Value2 = 2
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '0'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = '***.*'
$ws.Range("F22").Value2 = 4
$ws.Range("G22").Value2 = 3
$ws.Range("H22").Value2 = 33.333333333333
$ws.Range("I22").Value2 = 44
$ws.Range("K22").Value2 = 4.761904761904
$ws.Range("L22").Value2 = 4.761904761904
$ws.Range("M22").Value2 = -13.725490196078

# Row 23 (Housing)
$ws.Range("C23").Value2 = 6
$ws.Range("D23").Value2 = 6
$ws.Range("E23").Value2 = 0
$ws.Range("G23").Value2 = 27
$ws.Range("H23").Value2 = -29.629629629629
$ws.Range("I23").Value2 = 305
$ws.Range("J23").Value2 = 347
$ws.Range("K23").Value2 = -12.103746397694
$ws.Range("L23").Value2 = 0.993377483443
$ws.Range("M23").Value2 = 33.771929824561

# Row 24 (Petit Larceny)
$ws.Range("D24").Value2 = 32
$ws.Range("E24").Value2 = -53.125
$ws.Range("F24").Value2 = 77
$ws.Range("G24").Value2 = 97
$ws.Range("H24").Value2 = -20.618556701030
$ws.Range("I24").Value2 = 950
$ws.Range("J24").Value2 = 1022
$ws.Range("K24").Value2 = -7.045009784735
$ws.Range("L24").Value2 = 33.991537376586
$ws.Range("M24").Value2 = 8.200455580865

# Row 25 (Misd. Assault)
$ws.Range("C25").Value2 = 20
$ws.Range("D25").Value2 = 16
$ws.Range("E25").Value2 = 25
$ws.Range("F25").Value2 = 64
$ws.Range("G25").Value2 = 73
$ws.Range("H25").Value2 = -12.328767123287
$ws.Range("I25").Value2 = 866
$ws.Range("J25").Value2 = 850
$ws.Range("K25").Value2 = 1.882352941176
$ws.Range("L25").Value2 = 15.775401069518
$ws.Range("M25").Value2 = -25.728987993138

# Row 26 (UCR Rape*)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value2 = '0'
$ws.Range("E26").Value2 = -100
$ws.Range("F26").Value2 = 2
$ws.Range("G26").Value2 = 9
$ws.Range("H26").Value2 = -77.777777777777
$ws.Range("J26").Value2 = 52
$ws.Range("K26").Value2 = -25

# Row 27 (Other Sex Crimes)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value2 = '0'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '0'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = '***.*'
$ws.Range("G27").Value2 = 7
$ws.Range("H27").Value2 = -28.571428571428
$ws.Range("I27").Value2 = 74
$ws.Range("K27").Value2 = -3.896103896103
$ws.Range("L27").Value2 = -6.329113924050

# Row 28 (Shooting Vic.)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value2 = '0'
$ws.Range("F28").Value2 = 5
$ws.Range("G28").Value2 = 3
$ws.Range("H28").Value2 = 66.666666666666
$ws.Range("I28").Value2 = 43
$ws.Range("K28").Value2 = -41.891891891891
$ws.Range("L28").Value2 = -46.913580246913
$ws.Range("M28").Value2 = -55.670103092783
$ws.Range("N28").Value2 = -85.947712418300

# Row 29 (Shooting Inc.)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value2 = '0'
$ws.Range("F29").Value2 = 5
$ws.Range("G29").Value2 = 2
$ws.Range("H29").Value2 = 150
$ws.Range("I29").Value2 = 39
$ws.Range("K29").Value2 = -32.758620689655
$ws.Range("L29").Value2 = -42.647058823529
$ws.Range("M29").Value2 = -50.632911392405
$ws.Range("N29").Value2 = -86.021505376344

# Row 30 (Hate Crimes)
$ws.Range("D30").Value2 = 2
$ws.Range("E30").Value2 = -100
$ws.Range("G30").Value2 = 2
$ws.Range("H30").Value2 = -100
$ws.Range("J30").Value2 = 5
$ws.Range("K30").Value2 = -80
